# Add data for 2022-05-28:
#  - Sheet name "Through 2022-05-19" -> "Through 2022-05-20"
#  - Header I1 "2022 (through 05-19)" -> "2022 (through 05-20)"
#  - June 2022 count (I6) 65 -> 69
#  - Total 2022 count (I14) 617 -> 621

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-05-20"

$ws.Range("I1").Value = "2022 (through 05-20)"
$ws.Range("I6").Value = 69
$ws.Range("I14").Value = 621
